$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Range("H16").Value = 4000
$ws.Range("I16").Value = 4000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3770
$ws.Range("N16").Value = $null

# Row 17
$ws.Range("H17").Value = 296.8421
$ws.Range("J17").Value = 296.8421
$ws.Range("L17").Value = 890.5263
$ws.Range("N17").Value = -1226.5263

# Row 33
$ws.Range("H33").Value = 32146.594
$ws.Range("I33").Value = 40237.6
$ws.Range("J33").Value = 3250.1428
$ws.Range("K33").Value = 40237.6
$ws.Range("L33").Value = 3250.1428
$ws.Range("M33").Value = -40008.6
$ws.Range("N33").Value = -3708.1428

# Row 40
$ws.Range("H40").Value = 32203.516
$ws.Range("I40").Value = 47200.863
$ws.Range("J40").Value = 2208.818
$ws.Range("K40").Value = 47200.863
$ws.Range("L40").Value = 2208.818
$ws.Range("M40").Value = -47025.863
$ws.Range("N40").Value = -2558.818

# Row 51
$ws.Range("H51").Value = 9122.177
$ws.Range("I51").Value = 23436.2
$ws.Range("J51").Value = 3158
$ws.Range("K51").Value = 23436.2
$ws.Range("L51").Value = 3158
$ws.Range("M51").Value = -22952.2
$ws.Range("N51").Value = -4126

# Row 101
$ws.Range("H101").Value = 566
$ws.Range("I101").Value = 566
$ws.Range("K101").Value = 1698
$ws.Range("M101").Value = -76

# Row 121
$ws.Range("H121").Value = 1976.6666
$ws.Range("J121").Value = 1972
$ws.Range("L121").Value = 5916
$ws.Range("N121").Value = -9410

# Row 129
$ws.Range("H129").Value = 7828.933
$ws.Range("J129").Value = 1264.1538
$ws.Range("L129").Value = 3792.4614
$ws.Range("N129").Value = -13792.4614

# Row 137
$ws.Range("H137").Value = 1645.4231
$ws.Range("I137").Value = 1635.8422
$ws.Range("J137").Value = 1671.4286
$ws.Range("K137").Value = 4907.5266
$ws.Range("L137").Value = 5014.2858
$ws.Range("M137").Value = -2357.5266
$ws.Range("N137").Value = -10114.2858

# Row 141
$ws.Range("H141").Value = 3652.2222
$ws.Range("I141").Value = 3311.6667
$ws.Range("J141").Value = 4333.3335
$ws.Range("K141").Value = 9935.000100000001
$ws.Range("L141").Value = 13000.0005
$ws.Range("M141").Value = -4755.000100000001
$ws.Range("N141").Value = -23360.0005

$ws = $wb.Worksheets.Item("ARM")
# Row 43
$ws.Range("H43").Value = 6033
$ws.Range("J43").Value = 6033
$ws.Range("L43").Value = 6033
$ws.Range("N43").Value = -6659

# Row 61
$ws.Range("H61").Value = 1297.6897
$ws.Range("I61").Value = 1212.2693
$ws.Range("K61").Value = 1212.2693
$ws.Range("M61").Value = -1000.2693

# Row 70
$ws.Range("H70").Value = 45000
$ws.Range("J70").Value = 45000
$ws.Range("L70").Value = 45000
$ws.Range("N70").Value = -45540

# Row 73
$ws.Range("H73").Value = 45000
$ws.Range("J73").Value = 45000
$ws.Range("L73").Value = 45000
$ws.Range("N73").Value = -46872

# Row 95
$ws.Range("H95").Value = 28456.666
$ws.Range("J95").Value = 28456.666
$ws.Range("L95").Value = 28456.666
$ws.Range("N95").Value = -33948.666

# Row 136
$ws.Range("H136").Value = 1297.6897
$ws.Range("I136").Value = 1212.2693
$ws.Range("K136").Value = 3636.8079
$ws.Range("M136").Value = -1086.8079

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 287853
$ws.Range("I105").Value = 252195
$ws.Range("J105").Value = 335397
$ws.Range("K105").Value = 252195
$ws.Range("L105").Value = 335397
$ws.Range("M105").Value = -250448
$ws.Range("N105").Value = -338891

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 24617.69
$ws.Range("I31").Value = 49110.57
$ws.Range("J31").Value = 3186.4167
$ws.Range("K31").Value = 49110.57
$ws.Range("L31").Value = 3186.4167
$ws.Range("M31").Value = -48815.57
$ws.Range("N31").Value = -3776.4167

# Row 34
$ws.Range("H34").Value = 24617.69
$ws.Range("I34").Value = 49110.57
$ws.Range("J34").Value = 3186.4167
$ws.Range("K34").Value = 49110.57
$ws.Range("L34").Value = 3186.4167
$ws.Range("M34").Value = -48908.57
$ws.Range("N34").Value = -3590.4167

# Row 132
$ws.Range("H132").Value = 36588376
$ws.Range("J132").Value = 35717880
$ws.Range("L132").Value = 107153640
$ws.Range("N132").Value = -107158700

$ws = $wb.Worksheets.Item("CUL")
# Row 60
$ws.Range("H60").Value = 622
$ws.Range("J60").Value = 900
$ws.Range("L60").Value = 2700
$ws.Range("N60").Value = -3202

# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null

# Row 98
$ws.Range("H98").Value = 44883.08
$ws.Range("I98").Value = 801
$ws.Range("J98").Value = 55903.6
$ws.Range("K98").Value = 2403
$ws.Range("L98").Value = 167710.8
$ws.Range("M98").Value = -905
$ws.Range("N98").Value = -170706.8

# Row 122
$ws.Range("H122").Value = 4225.3335
$ws.Range("I122").Value = 381.22726
$ws.Range("J122").Value = 21139.4
$ws.Range("K122").Value = 3431.04534
$ws.Range("L122").Value = 190254.6
$ws.Range("M122").Value = -981.0453400000001
$ws.Range("N122").Value = -195154.6

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 58825652
$ws.Range("I97").Value = 66668960
$ws.Range("J97").Value = 840
$ws.Range("K97").Value = 66668960
$ws.Range("L97").Value = 840
$ws.Range("M97").Value = -66668464
$ws.Range("N97").Value = -1832

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1687576.6
$ws.Range("I46").Value = 300
$ws.Range("J46").Value = 2025032
$ws.Range("K46").Value = 300
$ws.Range("L46").Value = 2025032
$ws.Range("M46").Value = -112
$ws.Range("N46").Value = -2025408

# Row 48
$ws.Range("H48").Value = 8845.5
$ws.Range("I48").Value = 5860.6665
$ws.Range("J48").Value = 17800
$ws.Range("K48").Value = 5860.6665
$ws.Range("L48").Value = 17800
$ws.Range("M48").Value = -5199.6665
$ws.Range("N48").Value = -19122

# Row 136
$ws.Range("H136").Value = 2700
$ws.Range("I136").Value = 4100
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 12300
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -9750
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
# Row 24
$ws.Range("H24").Value = 25000
$ws.Range("J24").Value = 25000
$ws.Range("L24").Value = 25000
$ws.Range("N24").Value = -25460

# Row 39
$ws.Range("H39").Value = 13706.25
$ws.Range("I39").Value = 5625
$ws.Range("J39").Value = 16400
$ws.Range("K39").Value = 5625
$ws.Range("L39").Value = 16400
$ws.Range("M39").Value = -5212
$ws.Range("N39").Value = -17226

# Row 42
$ws.Range("H42").Value = 21949.889
$ws.Range("I42").Value = 13333.333
$ws.Range("J42").Value = 26258.166
$ws.Range("K42").Value = 13333.333
$ws.Range("L42").Value = 26258.166
$ws.Range("M42").Value = -12955.333
$ws.Range("N42").Value = -27014.166

# Row 43
$ws.Range("H43").Value = 52527.5
$ws.Range("I43").Value = 52527.5
$ws.Range("K43").Value = 52527.5
$ws.Range("M43").Value = -52378.5
